$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Styles-table bookkeeping artifact
# ---------------------------------------------------------------------------
# This workbook has accumulated an extra, unused "Arial 11" font entry in
# xl/styles.xml every time it picked up a fresh data drop (fonts count keeps
# climbing: 8 -> 9 -> 10 ...), without ever being referenced by a cell style.
# Reproduce the same bookkeeping side effect here: touch the workbook's
# built-in "Normal" style font and set it right back. The round trip leaves
# one orphaned font definition behind (fonts count 9 -> 10) while the
# "Normal" style / xf(0) keep pointing at the original font, so no visible
# cell formatting changes.
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Arial"
$normalStyle.Font.Name = "HP Simplified"

# ---------------------------------------------------------------------------
# Data edits
# ---------------------------------------------------------------------------
# Column AX ("Previous Doc"), rows 2-31: 4001967730 -> 0465572171
# Column K  ("Ship To Customer Name"), rows 2-31: dvdpdy341262 -> fbqdfs851092
#
# AX's new value is a numeric-looking string (leading zero). Assigning it
# through .Value directly gets auto-coerced to a Double (losing the leading
# zero and flipping the cell onto a brand-new quote-prefixed style). Instead,
# stage the literal text as a formula result on a scratch cell far outside
# the sheet's used range, then Copy / PasteSpecial (values only) it into each
# target cell -- this preserves the exact text, keeps the shared-string cell
# type, and leaves the cell's existing style untouched. The scratch cell is
# cleared (not just ClearContents) afterwards so the sheet's used range/
# dimension is unaffected.
$scratch = $ws.Cells.Item(1, 100)
$scratch.Formula = '="0465572171"'
$scratch.Copy()
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 50).PasteSpecial(-4163)
}
$scratch.Clear()

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 11).Value = "fbqdfs851092"
}
